# "generating monthly Report now Working ~Linus"
#
# Updates the Wareneingang (goods-receipt) sheet for the monthly report:
#  - fills in the missing "Zielland" (destination country) for row 2
#  - corrects the "Gewicht" (weight) totals for rows 2 and 3
#  - removes the duplicate/placeholder 4th row that was left over from testing
#  - leaves the selection/active view on the next empty row (E5), ready for
#    the next month's entries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: destination country was left blank - fill it in, and the weight
# needs to reflect the real shipment total.
$ws.Range("C2").Value = "DE"
$ws.Range("D2").Value = 1500

# Row 3: weight total correction.
$ws.Range("D3").Value = 1000

# Row 4 was a stray duplicate entry (same data as row 2, with no Zielland) -
# drop it now that the report only needs the two real rows.
$ws.Rows("4:4").Delete()

# Leave the view parked just past the last data row, matching where the
# next entry will be typed in.
$ws.Range("E5").Select()

# Restore the window geometry Excel normally persists for this workbook.
$win = $wb.Windows.Item(1)
$win.Left = 0
$win.Top = 760
$win.Width = 30240
$win.Height = 18880
